# Apply updated cryptocurrency price/volume figures scraped on
# Sun Oct  8 11:44:27 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to stay plain TEXT.
# The "Price" column holds numeric-looking strings (e.g. "23.20", "1.00",
# "27.817.62") that must be preserved exactly as text (trailing zeros and
# all) rather than being auto-coerced into numbers by Excel's normal entry
# parsing. A leading apostrophe is the standard Excel "treat as text" quote
# prefix; resetting the Style back to "Normal" afterwards drops the
# leftover quote-prefix formatting flag so the cell's style index is left
# exactly as it was.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

# --- Update Price (column D) and Volume(1h) (column E) figures for rows 2-43 ---
Set-TextValue 'D2' '27.817.62'
$ws.Range('E2').Value = '  -0.63%  '
Set-TextValue 'D3' '1.616.94'
$ws.Range('E3').Value = '  -1.50%  '
$ws.Range('E4').Value = '  -0.14%  '
Set-TextValue 'D5' '209.68'
$ws.Range('E6').Value = '  -0.57%  '
Set-TextValue 'D7' '0.999'
$ws.Range('E7').Value = '  -0.18%  '
Set-TextValue 'D8' '23.20'
$ws.Range('E8').Value = '  -1.38%  '
$ws.Range('E9').Value = '  -1.03%  '
Set-TextValue 'D10' '0.0609'
$ws.Range('E10').Value = '  -0.93%  '
Set-TextValue 'D11' '0.0877'
$ws.Range('E11').Value = '  -0.64%  '
Set-TextValue 'D12' '1.848.23'
$ws.Range('E12').Value = '  -1.40%  '
Set-TextValue 'D13' '1.616.20'
$ws.Range('E13').Value = '  -1.65%  '
$ws.Range('E14').Value = '  -2.11%  '
Set-TextValue 'D15' '0.557'
$ws.Range('E15').Value = '  -3.10%  '
Set-TextValue 'D16' '65.01'
$ws.Range('E16').Value = '  -1.31%  '
Set-TextValue 'D17' '27.801.68'
$ws.Range('E17').Value = '  -0.70%  '
Set-TextValue 'D18' '228.41'
$ws.Range('E18').Value = '  -3.35%  '
$ws.Range('E19').Value = '  -0.91%  '
Set-TextValue 'D20' '7.58'
$ws.Range('E20').Value = '  -0.53%  '
Set-TextValue 'D21' '0.999'
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('E22').Value = '  -2.04%  '
$ws.Range('E23').Value = '  -5.66%  '
$ws.Range('E24').Value = '  -3.64%  '
Set-TextValue 'D25' '154.47'
$ws.Range('E25').Value = '  +1.79%  '
$ws.Range('E26').Value = '  -1.31%  '
$ws.Range('E27').Value = '  -0.38%  '
Set-TextValue 'D28' '15.43'
$ws.Range('E28').Value = '  -1.67%  '
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('E30').Value = '  -1.58%  '
$ws.Range('E31').Value = '  -1.10%  '
Set-TextValue 'D32' '3.41'
$ws.Range('E32').Value = '  +1.99%  '
$ws.Range('E33').Value = '  -1.91%  '
Set-TextValue 'D34' '1.389.42'
$ws.Range('E34').Value = '  -1.97%  '
$ws.Range('E35').Value = '  -0.79%  '
Set-TextValue 'D36' '0.991'
$ws.Range('E36').Value = '  +10.17%  '
$ws.Range('E37').Value = '  -1.12%  '
$ws.Range('E38').Value = '  -0.44%  '
Set-TextValue 'D39' '0.552'
$ws.Range('E39').Value = '  -1.07%  '
Set-TextValue 'D40' '0.842'
$ws.Range('E40').Value = '  -4.65%  '
$ws.Range('E41').Value = '  -0.24%  '
Set-TextValue 'D42' '1.00'
$ws.Range('E42').Value = '  -2.41%  '
$ws.Range('E43').Value = '  -2.19%  '

# --- Rows 44 and 45 swapped order (FraxShare now ranks above Aave) ---
# Row 44 becomes FraxShare, row 45 becomes Aave, with refreshed price/volume data.
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D44' '5.44'
$ws.Range('E44').Value = '  -1.37%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D45' '65.35'
$ws.Range('E45').Value = '  -1.92%  '

# --- Update Price (column D) and Volume(1h) (column E) figures for rows 46-51 ---
Set-TextValue 'D46' '1.758.97'
$ws.Range('E46').Value = '  -1.38%  '
$ws.Range('E47').Value = '  -2.54%  '
Set-TextValue 'D48' '87.59'
$ws.Range('E48').Value = '  -0.36%  '
Set-TextValue 'D49' '0.0₆0102'
$ws.Range('E49').Value = '  -3.13%  '
$ws.Range('E50').Value = '  +0.29%  '
$ws.Range('E51').Value = '  -0.78%  '
